$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the placeholder "127.0.0.1" IP values with real per-server addresses.
$ws.Range("F2").Value = "192.168.1.113"
$ws.Range("F3").Value = "192.168.1.114"
$ws.Range("F4").Value = "192.168.1.115"
$ws.Range("F5").Value = "192.168.1.116"
$ws.Range("F6").Value = "192.168.1.117"

# Row 7 was a stray duplicate of the GameServer_1 entry (causing the
# "error in linux" due to a duplicate port/IP); clear its contents while
# keeping the existing cell formatting.
$ws.Range("A7:H7").ClearContents()

# Widen column F so the longer IP addresses are fully visible (best fit).
$ws.Columns.Item(6).ColumnWidth = 14.285714285714286

# Move the active selection to the now-empty row so it's ready for new data.
$ws.Range("A7:XFD7").Select()
